$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

$ws.Range("D2").Value = 3799.07
$ws.Range("E2").Value = -3799.07

$ws.Range("D4").Value = 4120.92
$ws.Range("E4").Value = 13379.08
$ws.Range("F4").Value = 0.2354811428571429
